$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sender's password (C2) was not correct; fix it by updating the value.
$ws.Range("C2").Value = "XXXXXX"
